$wb = $excel.ActiveWorkbook

# --- Operator sheet: update team composition inputs ---
# C10 (Type 1 count) and D10 (Type 2 count) both drop from 20 to 10.
# Dependent formulas (B10, B11, B4, B6) recalc automatically.
$wsOperator = $wb.Worksheets.Item("Operator")
$wsOperator.Range("C10").Value = 10
$wsOperator.Range("D10").Value = 10

# --- Aircraft sheet: update Range (mi) for Type 1 aircraft ---
$wsAircraft = $wb.Worksheets.Item("Aircraft")
$wsAircraft.Range("C2").Value = 150

# --- View/selection state ---
# Ports was the active tab before; move its selection without leaving it active.
$wsPorts = $wb.Worksheets.Item("Ports")
$wsPorts.Range("H23").Select()

# Aircraft becomes the active tab, with B9 selected.
$wsAircraft.Activate()
$wsAircraft.Range("B9").Select()
